$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell G6 held the text value "1,9" (Russian decimal comma) which Excel
# stored as a shared string (text number format). Fix the type error by
# switching the cell to General format and writing a real number (1.9)
# instead of text.
$ws.Range("G6").NumberFormat = "General"
$ws.Range("G6").Value = 1.9

# Move the active selection as recorded after the edit.
$ws.Range("B9").Select()
